# update code tinh luong
# - Đơn phụ phẫu 1: add a new completed order row (shifts the "Tổng" summary row down)
# - Rename the old "Lương" sheet to "Đơn phụ phẫu 2" and fill it with a new order row
# - Add a brand-new "Lương" sheet at the end with the recomputed payroll breakdown

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "Đơn phụ phẫu 1" — insert new row 3 (pushes the old Tổng row to 4)
# -----------------------------------------------------------------
$wsPP1 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

$wsPP1.Rows.Item(3).Insert()

$wsPP1.Cells.Item(3, 1).Value = "HD-LUXURY"
$wsPP1.Cells.Item(3, 2).Value = 596
$wsPP1.Cells.Item(3, 3).NumberFormat = "@"
$wsPP1.Cells.Item(3, 3).Value = "07-28-2024"
$wsPP1.Cells.Item(3, 4).Value = "SÓC TRĂNG"
$wsPP1.Cells.Item(3, 5).Value = "bảo trân"
$wsPP1.Cells.Item(3, 6).Value = "Cá nhân"
$wsPP1.Cells.Item(3, 7).Value = "Cắt mí"
$wsPP1.Cells.Item(3, 8).Value = "Trần Khánh Hiệp"
$wsPP1.Cells.Item(3, 9).Value = 50000

# update the (now shifted) Tổng row — row 4
$wsPP1.Cells.Item(4, 2).Value = 2
$wsPP1.Cells.Item(4, 9).Value = 50000

# -----------------------------------------------------------------
# 2) Reuse the old "Lương" sheet as "Đơn phụ phẫu 2" with fresh data
# -----------------------------------------------------------------
$wsPP2 = $wb.Worksheets.Item("Lương")
$wsPP2.Name = "Đơn phụ phẫu 2"
$wsPP2.Cells.Clear()

$pp2Headers = @("Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng", "Nguồn khách", "Tên dịch vụ", "Phụ phẫu 2", "Công phụ phẫu 2")
for ($c = 0; $c -lt $pp2Headers.Length; $c++) {
    $wsPP2.Cells.Item(1, $c + 1).Value = $pp2Headers[$c]
}

$wsPP2.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsPP2.Cells.Item(2, 2).Value = 583
$wsPP2.Cells.Item(2, 3).NumberFormat = "@"
$wsPP2.Cells.Item(2, 3).Value = "07-23-2024"
$wsPP2.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$wsPP2.Cells.Item(2, 5).Value = "nguyễn thị mỹ trinh"
$wsPP2.Cells.Item(2, 6).Value = "Cá nhân"
$wsPP2.Cells.Item(2, 7).Value = "Nâng mũi"
$wsPP2.Cells.Item(2, 8).Value = "Trần Khánh Hiệp"
$wsPP2.Cells.Item(2, 9).Value = 50000

$wsPP2.Cells.Item(3, 1).Value = "Tổng"
$wsPP2.Cells.Item(3, 2).Value = 1
$wsPP2.Cells.Item(3, 9).Value = 50000

# -----------------------------------------------------------------
# 3) Add a new "Lương" sheet (after "Đơn phụ phẫu 2") with the
#    recalculated payroll summary
# -----------------------------------------------------------------
$wsLuong = $wb.Worksheets.Add($null, $wsPP2)
$wsLuong.Name = "Lương"

$luongRows = @(
    @("Danh mục lương", 10),
    @("Tổng công tại CẦN THƠ", 0),
    @("Phụ cấp tại CẦN THƠ", 0),
    @("Lương công tác tại CẦN THƠ", 0),
    @("Lương cơ bản tại CẦN THƠ", $null),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Ứng lương tại CẦN THƠ", -0),
    @("Tổng công tại LONG XUYÊN", 0),
    @("Phụ cấp tại LONG XUYÊN", 0),
    @("Lương công tác tại LONG XUYÊN", 0),
    @("Lương cơ bản tại LONG XUYÊN", $null),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Ứng lương tại LONG XUYÊN", -0),
    @("Tổng công tại SÓC TRĂNG", 26),
    @("Phụ cấp tại SÓC TRĂNG", 910000),
    @("Lương cơ bản tại SÓC TRĂNG", 3825714.285714285),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 50000),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 50000),
    @("Ứng lương tại SÓC TRĂNG", -2000000),
    @("Tổng lương tại CẦN THƠ", 0),
    @("Tổng lương tại LONG XUYÊN", 0),
    @("Tổng lương tại SÓC TRĂNG", 2835714.285714285),
    @("Tổng lương", 2835714.285714285)
)

for ($r = 0; $r -lt $luongRows.Length; $r++) {
    $row = $luongRows[$r]
    $wsLuong.Cells.Item($r + 1, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $wsLuong.Cells.Item($r + 1, 2).Value = $row[1]
    }
}
